$wb = $excel.ActiveWorkbook

$home = $wb.Worksheets.Item("Home")
$checking = $wb.Worksheets.Item("Checking")
$savings = $wb.Worksheets.Item("Savings")
$business = $wb.Worksheets.Item("Business")

# New log entry (transaction) added to the Home sheet as row 15.
# The Date column holds plain text like "2021/10/01" (not a real date), so
# copy an existing date cell's text over via PasteSpecial (values only) to
# avoid Excel auto-converting the literal into a date serial number.
$home.Cells.Item(5, 6).Copy()
$home.Cells.Item(15, 6).PasteSpecial(-4163)
$home.Cells.Item(15, 7).Value = "16:32:37"
$home.Cells.Item(15, 8).Value = "Checkings"
$home.Cells.Item(15, 9).Value = 123
$home.Cells.Item(15, 10).Value = "test123"
$home.Cells.Item(15, 11).Value = 1376

# Mirror the same entry onto the Checking sheet (row 8), which only tracks
# "Checkings" category rows
$checking.Cells.Item(3, 6).Copy()
$checking.Cells.Item(8, 6).PasteSpecial(-4163)
$checking.Cells.Item(8, 7).Value = "16:32:37"
$checking.Cells.Item(8, 8).Value = "Checkings"
$checking.Cells.Item(8, 9).Value = 123
$checking.Cells.Item(8, 10).Value = "test123"
$checking.Cells.Item(8, 11).Value = 1376

# Replace the cross-sheet "New Total" link formulas with their current
# static values
$checking.Range("N2").Value = 1499
$savings.Range("N2").Value = 27300
$business.Range("N2").Value = 315000111

# Update selections: user ends up on the Business sheet
$home.Range("P2").Select()
$checking.Range("N2").Select()
$savings.Range("N2").Select()
$business.Range("N2").Select()
$business.Activate()
